$d = $word.ActiveDocument

$replacements = @(
    @("64÷9=", "49÷9="),
    @("65÷9=", "50÷6="),
    @("56÷5=", "18÷7="),
    @("76÷8=", "88÷9="),
    @("39÷2=", "39÷6="),
    @("24÷3=", "27÷9="),
    @("81÷7=", "62÷8="),
    @("58÷7=", "59÷3="),
    @("46÷5=", "44÷7="),
    @("87÷8=", "43÷7="),
    @("69÷5=", "40÷6="),
    @("75÷8=", "22÷5="),
    @("53÷8=", "64÷4="),
    @("73÷7=", "68÷2="),
    @("93÷6=", "16÷5="),
    @("59÷5=", "66÷5="),
    @("99÷8=", "30÷7="),
    @("90÷5=", "83÷3="),
    @("86÷5=", "96÷9="),
    @("29÷4=", "47÷5="),
    @("92÷4=", "38÷5="),
    @("33÷9=", "93÷7="),
    @("46÷3=", "56÷3="),
    @("47÷9=", "67÷7="),
    @("21÷2=", "55÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
